# Kakao LogUp / LogIn related seed-data update
# - ErrorMessages sheet gains a new row describing "modulse/user.User.logIn" -> "deleted"
# - selection/active-cell bookkeeping is refreshed on the Status and ErrorMessages sheets

$wb = $excel.ActiveWorkbook

$wsStatus = $wb.Worksheets.Item(1)        # "Status" sheet
$wsErrors = $wb.Worksheets.Item(2)        # "ErrorMessages" sheet

# ---------------------------------------------------------------------------
# ErrorMessages: append row 32 with the new log-in error entry
# ---------------------------------------------------------------------------

# Copy the look & feel (style) of the row directly above before filling values
# so the new cells keep the same formatting (s="1") as the rest of the table.
$wsErrors.Range("B31").Copy($wsErrors.Range("B32"))
$wsErrors.Range("D31").Copy($wsErrors.Range("D32"))

$wsErrors.Range("A32").Value = 30
$wsErrors.Range("B32").Value = "modulse/user.User.logIn"
$wsErrors.Range("C32").Value = 400
$wsErrors.Range("D32").Value = "deleted"

# ---------------------------------------------------------------------------
# Refresh view state (active cell / selection) on both sheets
# ---------------------------------------------------------------------------

$wsStatus.Activate()
$wsStatus.Range("B4").Select()

$wsErrors.Activate()
$wsErrors.Range("D32").Select()
